$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New players entered into the list, in the order they were typed in
$ws.Range("A41").Value = "MonkeyDisco"
$ws.Range("A53").Value = "eXeQtr"
$ws.Range("A54").Value = "KissU"

# Row 23: placeholder "DEADZONE" name is replaced with the real player name
$ws.Range("A23").Value = "VINDICATER"

$ws.Range("A55").Value = "dodiz"
$ws.Range("A56").Value = "Lyne"
$ws.Range("A57").Value = "KwEne"
$ws.Range("A58").Value = "Zarich"
$ws.Range("A59").Value = "DonaldDuck"
$ws.Range("A60").Value = "Dewestator"
$ws.Range("A61").Value = "Carbon"

# Restore the view/selection state recorded in the saved workbook
$ws.Range("E64").Select()
$excel.ActiveWindow.ScrollRow = 38
